$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4335.268898379608
$ws.Range("C3").Value = 4276.118057490181
$ws.Range("C4").Value = 4276.118057490181
$ws.Range("C5").Value = 4114.02442874095
$ws.Range("C6").Value = 4088.449422597257
$ws.Range("C7").Value = 4088.449422597257
$ws.Range("C8").Value = 4088.449422597257
$ws.Range("C9").Value = 4088.449422597257
$ws.Range("C10").Value = 4088.449422597257
$ws.Range("C11").Value = 4023.44003937617
$ws.Range("C12").Value = 4023.44003937617
